# Add "continue" and "new game" strings/rows to both language sheets
# (en = sheet1, es = sheet2), matching the source workbook's layout.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "en"
$ws2 = $wb.Worksheets.Item(2)   # "es"

# --- Key column (A) values, identical on both sheets ---------------------
$ws1.Range("A55").Value = "new_game"
$ws2.Range("A55").Value = "new_game"

$ws1.Range("A56").Value = "continue"
$ws2.Range("A56").Value = "continue"

# --- Translated values (B) -------------------------------------------------
$ws2.Range("B56").Value = "CONTINUAR"
$ws1.Range("B56").Value = "CONTINUE"

$ws2.Range("B55").Value = "NUEVO JUEGO"
$ws1.Range("B55").Value = "NEW GAME"

# The "en" sheet uses the vertically-centered style (style index 2, as used
# by the other rows near the bottom of that sheet) for column B.
$ws1.Range("B55").VerticalAlignment = -4108
$ws1.Range("B56").VerticalAlignment = -4108

# --- Selection / active sheet state ---------------------------------------
# Select B55 on the "es" sheet first, then on the "en" sheet last so that
# "en" ends up as the active sheet/tab, matching the target workbook.
$null = $ws2.Range("B55").Select()
$null = $ws1.Range("B55").Select()
$null = $ws1.Activate()

Write-Output "done"
